$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 128; everything from 128..216 shifts down to 129..217,
# and the sheet dimension grows from T216 to T217.
$ws.Rows("128:128").Insert()

# Populate the newly inserted row 128 with the new record (a daily Ciruela price
# observation for Macroferia Regional de Talca).
$ws.Cells.Item(128, 1).Value  = 5
$ws.Cells.Item(128, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(128, 3).Value  = "Maule"
$ws.Cells.Item(128, 4).Value  = 45040
$ws.Cells.Item(128, 5).Value  = 7
$ws.Cells.Item(128, 6).Value  = "Fruta"
$ws.Cells.Item(128, 7).Value  = 100103
$ws.Cells.Item(128, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(128, 9).Value  = 100103002
$ws.Cells.Item(128, 10).Value = "Ciruela"
$ws.Cells.Item(128, 11).Value = "Angeleno"
$ws.Cells.Item(128, 12).Value = "Primera"
$ws.Cells.Item(128, 13).Value = 220
$ws.Cells.Item(128, 14).Value = 10000
$ws.Cells.Item(128, 15).Value = 10000
$ws.Cells.Item(128, 16).Value = 10000
$ws.Cells.Item(128, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(128, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(128, 19).Value = 556
$ws.Cells.Item(128, 20).Value = 18
